$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B6").Value = "Shakib Ahmed"
$ws.Range("C6").Value = "Shak Forid"
$ws.Range("D6").Value = "Aklima"
$ws.Range("E6").Value = 677287
$ws.Range("G6").Value = "Computer Science & Technology (85)"
$ws.Range("H6").Value = "23/06/2003"
